# Updated cryptos list -- refresh Price (column D) and Volume(1h) (column E)
# values for the crypto table on Sheet1, rows 2-51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '24.710.12'; E = '  +0.88%  ' },
    @{ Row = 3; D = '1.693.63'; E = '  +0.78%  ' },
    @{ Row = 4; D = '1.006'; E = '  +0.46%  ' },
    @{ Row = 5; D = '310.21'; E = '  +1.42%  ' },
    @{ Row = 6; D = '1.000'; E = '  +0.38%  ' },
    @{ Row = 7; D = '0.3711'; E = '  +0.69%  ' },
    @{ Row = 8; D = '48.94'; E = '  +1.96%  ' },
    @{ Row = 9; D = '0.3406'; E = '  -0.50%  ' },
    @{ Row = 10; D = '1.201'; E = '  +3.35%  ' },
    @{ Row = 11; D = '0.07417'; E = '  +2.67%  ' },
    @{ Row = 12; D = '1.002'; E = '  +0.44%  ' },
    @{ Row = 13; D = '6.278'; E = '  +2.83%  ' },
    @{ Row = 14; D = '20.81'; E = '  +3.15%  ' },
    @{ Row = 15; D = '6.929'; E = '  +3.25%  ' },
    @{ Row = 16; D = '1.690.63'; E = '  +0.79%  ' },
    @{ Row = 17; D = '0.00001115'; E = '  +1.28%  ' },
    @{ Row = 18; D = '0.06698'; E = '  +0.75%  ' },
    @{ Row = 19; D = $null; E = '  +0.44%  ' },
    @{ Row = 20; D = '82.81'; E = '  +2.69%  ' },
    @{ Row = 21; D = '17.07'; E = '  +3.89%  ' },
    @{ Row = 22; D = '6.279'; E = '  +3.33%  ' },
    @{ Row = 23; D = '12.81'; E = '  +5.98%  ' },
    @{ Row = 24; D = '24.740.07'; E = '  +1.28%  ' },
    @{ Row = 25; D = $null; E = '  +1.01%  ' },
    @{ Row = 26; D = '2.743'; E = '  +3.40%  ' },
    @{ Row = 27; D = '20.09'; E = '  +3.47%  ' },
    @{ Row = 28; D = $null; E = '  -3.15%  ' },
    @{ Row = 29; D = '131.21'; E = '  +3.33%  ' },
    @{ Row = 30; D = '1.878.18'; E = '  +0.81%  ' },
    @{ Row = 31; D = '1.235'; E = '  +26.01%  ' },
    @{ Row = 32; D = '6.659'; E = '  +6.59%  ' },
    @{ Row = 33; D = '4.210'; E = '  +4.61%  ' },
    @{ Row = 34; D = '13.51'; E = '  +9.55%  ' },
    @{ Row = 35; D = '1.749'; E = '  +3.37%  ' },
    @{ Row = 36; D = '0.08661'; E = '  +2.93%  ' },
    @{ Row = 37; D = '5.489'; E = '  +3.23%  ' },
    @{ Row = 38; D = '0.06566'; E = '  +3.03%  ' },
    @{ Row = 39; D = '9.011'; E = '  +3.84%  ' },
    @{ Row = 40; D = '0.02387'; E = '  +3.38%  ' },
    @{ Row = 41; D = '0.2193'; E = '  +4.94%  ' },
    @{ Row = 42; D = '1.255'; E = '  +0.82%  ' },
    @{ Row = 43; D = '0.6349'; E = '  +4.28%  ' },
    @{ Row = 44; D = '1.001'; E = '  +0.41%  ' },
    @{ Row = 45; D = '13.63'; E = '  +4.99%  ' },
    @{ Row = 46; D = '3.804'; E = '  +1.19%  ' },
    @{ Row = 47; D = '0.6021'; E = '  +2.66%  ' },
    @{ Row = 48; D = '2.091'; E = '  +4.15%  ' },
    @{ Row = 49; D = '127.83'; E = '  +1.83%  ' },
    @{ Row = 50; D = $null; E = '  +0.82%  ' },
    @{ Row = 51; D = '78.74'; E = '  +4.08%  ' }
)

# Column D holds numeric-looking text (e.g. "1.006", "24.710.12") that must
# stay plain text, matching the source file's inlineStr cells, not become
# a real number. Prefixing with an apostrophe forces text entry; resetting
# the style back to Normal afterwards undoes the implicit "@" text format
# that gets stamped on the cell so no stray style is introduced.
foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($r, 4)
        $dCell.Value = "'" + $u.D
        $dCell.Style = "Normal"
    }
    $ws.Cells.Item($r, 5).Value = $u.E
}
